# Generate Report for Handback
# Update timestamps (and one status value) recorded on the handback status report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-17 04:15:19"
$wsOverview.Range("G5").Value = "2016-08-17 04:15:19"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-17 04:15:15"
$wsZhCn.Range("H5").Value = "2016-08-17 04:15:15"
$wsZhCn.Range("K4").Value = "2016-08-17 04:15:32"
$wsZhCn.Range("K5").Value = "2016-08-17 04:15:32"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-17 04:15:19"
$wsDeDe.Range("H5").Value = "2016-08-17 04:15:19"
$wsDeDe.Range("K4").Value = "2016-08-17 04:15:39"
$wsDeDe.Range("K5").Value = "2016-08-17 04:15:39"
